# edit.ps1 - reproduce the captured change:
#   1. Slide 6's table switches to the built-in table style
#      {1B8E2D73-0249-414C-9BEE-9F03BB3D26A4}
#   2. The presentation's theme (Slide Master) is changed from the
#      custom "Integral" palette to the standard "Office Theme" palette.

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 -------------------------------------------
$slide = $p.Slides.Item(6)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{1B8E2D73-0249-414C-9BEE-9F03BB3D26A4}")
    }
}

# --- 2. Re-colour the theme to the stock "Office Theme" palette ----------
function RGBVal($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

# index -> (r,g,b) for dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink
$officeColors = @(
    @(0x00,0x00,0x00),  # 1 dk1
    @(0xFF,0xFF,0xFF),  # 2 lt1
    @(0x44,0x54,0x6A),  # 3 dk2
    @(0xE7,0xE6,0xE6),  # 4 lt2
    @(0x5B,0x9B,0xD5),  # 5 accent1
    @(0xED,0x7D,0x31),  # 6 accent2
    @(0xA5,0xA5,0xA5),  # 7 accent3
    @(0xFF,0xC0,0x00),  # 8 accent4
    @(0x44,0x72,0xC4),  # 9 accent5
    @(0x70,0xAD,0x47),  # 10 accent6
    @(0x05,0x63,0xC1),  # 11 hlink
    @(0x95,0x4F,0x72)   # 12 folHlink
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $rgb = $officeColors[$i - 1]
    $themeColors.Colors($i).RGB = RGBVal $rgb[0] $rgb[1] $rgb[2]
}
